$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force-text cells (column D) use a leading apostrophe so Excel keeps them as
# plain text instead of auto-coercing number-looking strings (e.g. "210.19")
# into numeric cells; Style is reset to Normal afterwards so no stray
# quote-prefix formatting is left behind on the cell.

$ws.Range("D2").Value = "'28.304.46"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.49%  "
$ws.Range("D3").Value = "'1.552.39"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.38%  "
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").Value = "'210.13"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.54%  "
$ws.Range("E6").Value = "  -1.93%  "
$ws.Range("E7").Value = "  -0.18%  "
$ws.Range("D8").Value = "'23.80"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.39%  "
$ws.Range("E9").Value = "  -1.86%  "
$ws.Range("E10").Value = "  -1.72%  "
$ws.Range("E11").Value = "  +0.08%  "
$ws.Range("D12").Value = "'1.774.01"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.42%  "
$ws.Range("D13").Value = "'1.545.46"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.87%  "
$ws.Range("D14").Value = "'28.280.04"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.56%  "
$ws.Range("E15").Value = "  -1.59%  "
$ws.Range("E16").Value = "  -2.44%  "
$ws.Range("D17").Value = "'60.58"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.99%  "
$ws.Range("D18").Value = "'227.96"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.29%  "
$ws.Range("E19").Value = "  -0.84%  "
$ws.Range("D20").Value = "'0.0₃0675"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.67%  "
$ws.Range("E21").Value = "  -0.15%  "
$ws.Range("E22").Value = "  +0.72%  "
$ws.Range("E23").Value = "  -2.96%  "
$ws.Range("E24").Value = "  -2.98%  "
$ws.Range("D25").Value = "'151.24"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Value = "'14.76"
$ws.Range("D26").Style = "Normal"
$ws.Range("E27").Value = "  -1.63%  "
$ws.Range("E28").Value = "  -0.22%  "
$ws.Range("E31").Value = "  -4.30%  "
$ws.Range("E32").Value = "  -1.49%  "
$ws.Range("D33").Value = "'1.388.14"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.57%  "
$ws.Range("E34").Value = "  -3.01%  "
$ws.Range("E35").Value = "  +2.92%  "
$ws.Range("E36").Value = "  -4.35%  "
$ws.Range("E37").Value = "  -1.29%  "
$ws.Range("E40").Value = "  -2.58%  "
$ws.Range("E41").Value = "  +0.80%  "
$ws.Range("E42").Value = "  -0.22%  "
$ws.Range("D43").Value = "'0.778"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.00%  "
$ws.Range("E44").Value = "  -1.35%  "
$ws.Range("E45").Value = "  -2.60%  "
$ws.Range("D46").Value = "'61.83"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.22%  "
$ws.Range("D47").Value = "'1.686.77"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.49%  "
$ws.Range("D48").Value = "'0.905"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -5.82%  "
$ws.Range("D49").Value = "'85.71"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.07%  "
$ws.Range("D50").Value = "'42.56"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +6.83%  "
$ws.Range("E51").Value = "  -0.13%  "
